$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "subject"
$ws.Range("B1").Value = "condition"

$ws.Range("B1").Select()
